$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.747.26"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.76%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.882.21"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.28%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.001"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "237.33"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +1.86%  "
$ws.Range("E6").Value = "  +0.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4766"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.95%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2831"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +4.18%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06497"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.50%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.72"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +15.98%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.913.86"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +4.97%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07564"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +1.85%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "95.12"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +13.77%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.073"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +2.70%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.6498"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "297.61"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +31.58%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.734.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.96%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.11"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +6.39%  "
$ws.Range("E19").Value = "  +0.06%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007481"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.135.12"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +2.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "1.001"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.09%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.125"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +4.61%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.116"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.69%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "168.72"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.63%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.245"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("E27").Value = "  +9.96%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.950"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +4.80%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.1058"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +2.54%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.354"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.42%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.161"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.12%  "
$ws.Range("E32").Value = "  +3.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05002"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +4.05%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.169"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.56%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.7180"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +2.16%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.715"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01923"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +2.97%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.699"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.82%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.054"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +7.47%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.8963"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.75%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "106.86"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +2.35%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4176"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +4.25%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.575"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +1.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.09%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.322"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +4.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.1216"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.91%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.856"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.65%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "34.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +4.79%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05630"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +2.15%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.379"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +1.93%  "
